$wb = $excel.ActiveWorkbook

# --- Data correction on "raw benthic data": rows 42-62 in column A were
# mislabeled as site "PRWI-MARU" but should be "PRWI-MAWI". Correct the
# site name and bring the cell formatting in line with the rest of the row
# (which had already been formatted correctly).
$ws3 = $wb.Worksheets.Item("raw benthic data")
$ws3.Range("A42:A62").Value = "PRWI-MAWI"
for ($r = 42; $r -le 62; $r++) {
    $srcCell = $ws3.Cells.Item($r, 2)
    $dstCell = $ws3.Cells.Item($r, 1)
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# --- Update active sheet / selections to reflect the saved view state.
# "raw benthic data" keeps its own remembered selection (C11), no longer
# scrolled to A40 and no longer the tab shown when reopened.
$ws3.Activate() | Out-Null
$ws3.Range("C11").Select() | Out-Null

# "water chem" becomes the active tab with B2 selected. Activate it last so
# it is the tab shown when the workbook is reopened.
$ws1 = $wb.Worksheets.Item("water chem")
$ws1.Activate() | Out-Null
$ws1.Range("B2").Select() | Out-Null
